# Append the 2025-07-02 bitcoin buy record as a new row (row 28) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (matching the other recently
# appended rows, e.g. A10, A23-A27, which store "MM/DD/YYYY" strings
# rather than Excel date serials). Prefixing with a single quote forces
# Excel to treat the entry as text instead of auto-converting it to a
# date value; resetting the style back to Normal afterwards clears the
# "quote prefix" flag that operation leaves behind, so the cell ends up
# with the same default (unstyled) formatting as its neighbours.
$ws.Range("A28").Value = "'07/02/2025"
$ws.Range("A28").Style = "Normal"

$ws.Range("B28").Value = -0.01231811
$ws.Range("C28").Value = -4059.064255798982
$ws.Range("D28").Value = 50
